# The edit inserts "또 다른 " right before "문제점은" in the sentence that
# begins "(*) 그외 문제점은 Linear Regression의 가설은 ...".
# Doing so also relocates the document's (single) "_GoBack" bookmark from its
# old spot (in the following paragraph, right after "이") to the newly typed
# text, matching how Word stamps the last-edit bookmark at the insertion
# point.

$d = $word.ActiveDocument

# Find the insertion point: right after "그외 " (with its trailing space),
# immediately before "문제점은".
$found = $d.Content
$found.Find.Execute("그외 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $found.End

$insertRange = $d.Range($insertPoint, $insertPoint)
$insertRange.InsertBefore("또 다른 ")

# Re-locate the just-inserted text so we can (a) keep it in its own run
# (toggling a character property on/off forces the engine to stop merging it
# with neighboring identically-formatted runs, the same way Word keeps freshly
# typed text in a distinct run), and (b) find where the new "_GoBack" bookmark
# belongs.
$newText = $d.Content
$newText.Find.Execute("또 다른 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newText.Font.Bold = $true
$newText.Font.Bold = $false

# Re-adding a bookmark named "_GoBack" removes any previous one elsewhere in
# the document and places the only instance at the new location -- exactly
# Word's behavior of tracking the most recent edit point.
$bookmarkPos = $newText.End
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
